$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.252.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.372.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.372.13"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.47"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.45%  "
$ws.Range("E11").Value = "  +1.49%  "
$ws.Range("E12").Value = "  -0.60%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.946.99"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.10%  "
$ws.Range("E14").Value = "  +2.42%  "
$ws.Range("E15").Value = "  +1.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.00"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.375.07"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.372.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("E19").Value = "  +0.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "375.48"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.554"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.508.68"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000126"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "71.46"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.69"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.81%  "
$ws.Range("E29").Value = "  +0.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.45"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.33%  "
$ws.Range("E32").Value = "  +3.36%  "
$ws.Range("E33").Value = "  +1.54%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("E35").Value = "  +0.27%  "
$ws.Range("E36").Value = "  -5.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.82"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.85%  "
$ws.Range("E38").Value = "  -1.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "165.70"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0774"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.96%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("E42").Value = "  +1.74%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.40"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.29%  "
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.98%  "
$ws.Range("E48").Value = "  -1.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.59"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.348.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.60%  "
$ws.Range("E51").Value = "  +0.38%  "
